# Auto-generated Excel COM-interop script
# Applies scheduled-runner market-price refresh to columns H-N (currentAveragePrice* / LevePrice* / LeveProfit*)
# across the eight job sheets, matching the authoritative diff exactly.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 1898.25
$ws.Range("I43").Value = 1020.25
$ws.Range("J43").Value = 2337.25
$ws.Range("K43").Value = 1020.25
$ws.Range("L43").Value = 2337.25
$ws.Range("M43").Value = -951.25
$ws.Range("N43").Value = -2475.25
$ws.Range("H113").Value = 6317.2666
$ws.Range("J113").Value = 6896.923
$ws.Range("L113").Value = 6896.923
$ws.Range("N113").Value = -13404.923
$ws.Range("H138").Value = 3358.14
$ws.Range("I138").Value = 932.3333
$ws.Range("J138").Value = 4397.7715
$ws.Range("K138").Value = 2796.9999
$ws.Range("L138").Value = 13193.3145
$ws.Range("M138").Value = 2343.0001
$ws.Range("N138").Value = -23473.3145
$ws.Range("H141").Value = 9156.704
$ws.Range("I141").Value = 9777.125
$ws.Range("J141").Value = 4193.3335
$ws.Range("K141").Value = 29331.375
$ws.Range("L141").Value = 12580.0005
$ws.Range("M141").Value = -24151.375
$ws.Range("N141").Value = -22940.0005

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5727.1934
$ws.Range("I32").Value = 4262.702
$ws.Range("J32").Value = 10315.934
$ws.Range("K32").Value = 4262.702
$ws.Range("L32").Value = 10315.934
$ws.Range("M32").Value = -3975.702
$ws.Range("N32").Value = -10889.934
$ws.Range("H61").Value = 1163.88
$ws.Range("I61").Value = 1066.4166
$ws.Range("J61").Value = 1253.8462
$ws.Range("K61").Value = 1066.4166
$ws.Range("L61").Value = 1253.8462
$ws.Range("M61").Value = -854.4166
$ws.Range("N61").Value = -1677.8462
$ws.Range("H76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("H79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("H80").Value = 31516.666
$ws.Range("J80").Value = 31516.666
$ws.Range("L80").Value = 31516.666
$ws.Range("N80").Value = -33512.666
$ws.Range("H83").Value = 31516.666
$ws.Range("J83").Value = 31516.666
$ws.Range("L83").Value = 94549.99800000001
$ws.Range("N83").Value = -104533.998
$ws.Range("H132").Value = 1752.3208
$ws.Range("I132").Value = 1160.4318
$ws.Range("K132").Value = 3481.2954
$ws.Range("M132").Value = -951.2954
$ws.Range("H136").Value = 1163.88
$ws.Range("I136").Value = 1066.4166
$ws.Range("J136").Value = 1253.8462
$ws.Range("K136").Value = 3199.2498
$ws.Range("L136").Value = 3761.5386
$ws.Range("M136").Value = -649.2498000000001
$ws.Range("N136").Value = -8861.5386
$ws.Range("H137").Value = 40780
$ws.Range("J137").Value = 40780
$ws.Range("L137").Value = 40780
$ws.Range("N137").Value = -50980
$ws.Range("N76").ClearContents()
$ws.Range("N79").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H64").Value = 668.9
$ws.Range("I64").Value = 1505
$ws.Range("J64").Value = 459.875
$ws.Range("K64").Value = 1505
$ws.Range("L64").Value = 459.875
$ws.Range("M64").Value = -1280
$ws.Range("N64").Value = -909.875
$ws.Range("H67").Value = 668.9
$ws.Range("I67").Value = 1505
$ws.Range("J67").Value = 459.875
$ws.Range("K67").Value = 1505
$ws.Range("L67").Value = 459.875
$ws.Range("M67").Value = -725
$ws.Range("N67").Value = -2019.875
$ws.Range("H134").Value = 2326.0833
$ws.Range("I134").Value = 1329.7241
$ws.Range("J134").Value = 3846.842
$ws.Range("K134").Value = 3989.1723
$ws.Range("L134").Value = 11540.526
$ws.Range("M134").Value = -1454.1723
$ws.Range("N134").Value = -16610.526
$ws.Range("H137").Value = 43730
$ws.Range("J137").Value = 43730
$ws.Range("L137").Value = 43730
$ws.Range("N137").Value = -53930

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 3147.5557
$ws.Range("I122").Value = 2200
$ws.Range("J122").Value = 3266
$ws.Range("K122").Value = 6600
$ws.Range("L122").Value = 9798
$ws.Range("M122").Value = -4150
$ws.Range("N122").Value = -14698
$ws.Range("H132").Value = 1787.6177
$ws.Range("I132").Value = 992.63635
$ws.Range("K132").Value = 2977.90905
$ws.Range("M132").Value = -447.9090500000002
$ws.Range("H134").Value = 3258.6296
$ws.Range("I134").Value = 3762.697
$ws.Range("J134").Value = 2466.524
$ws.Range("K134").Value = 11288.091
$ws.Range("L134").Value = 7399.572
$ws.Range("M134").Value = -8753.091
$ws.Range("N134").Value = -12469.572

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1204.4894
$ws.Range("I5").Value = 426
$ws.Range("J5").Value = 1950.5416
$ws.Range("K5").Value = 1278
$ws.Range("L5").Value = 5851.6248
$ws.Range("M5").Value = -1166
$ws.Range("N5").Value = -6075.6248
$ws.Range("H131").Value = 7143620.5
$ws.Range("J131").Value = 851.36206
$ws.Range("L131").Value = 2554.08618
$ws.Range("N131").Value = -12634.08618
$ws.Range("H132").Value = 2034.7142
$ws.Range("I132").Value = 797.2
$ws.Range("J132").Value = 2722.2222
$ws.Range("K132").Value = 7174.8
$ws.Range("L132").Value = 24499.9998
$ws.Range("M132").Value = -4644.8
$ws.Range("N132").Value = -29559.9998
$ws.Range("H135").Value = 1204.4894
$ws.Range("I135").Value = 426
$ws.Range("J135").Value = 1950.5416
$ws.Range("K135").Value = 3834
$ws.Range("L135").Value = 17554.8744
$ws.Range("M135").Value = -1299
$ws.Range("N135").Value = -22624.8744
$ws.Range("H140").Value = 2832.3914
$ws.Range("I140").Value = 3813.182
$ws.Range("J140").Value = 1933.3334
$ws.Range("K140").Value = 11439.546
$ws.Range("L140").Value = 5800.0002
$ws.Range("M140").Value = -6259.545999999998
$ws.Range("N140").Value = -16160.0002

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H62").Value = 37950
$ws.Range("J62").Value = 37950
$ws.Range("L62").Value = 37950
$ws.Range("N62").Value = -39322
$ws.Range("H65").Value = 37950
$ws.Range("J65").Value = 37950
$ws.Range("L65").Value = 113850
$ws.Range("N65").Value = -120714
$ws.Range("H122").Value = 4870
$ws.Range("I122").Value = 2340
$ws.Range("J122").Value = 7400
$ws.Range("K122").Value = 7020
$ws.Range("L122").Value = 22200
$ws.Range("M122").Value = -4570
$ws.Range("N122").Value = -27100
$ws.Range("H126").Value = 2061.98
$ws.Range("I126").Value = 2061.98
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 6185.940000000001
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -3715.940000000001
$ws.Range("H132").Value = 2843.3872
$ws.Range("I132").Value = 1756.1578
$ws.Range("J132").Value = 4564.8335
$ws.Range("K132").Value = 5268.4734
$ws.Range("L132").Value = 13694.5005
$ws.Range("M132").Value = -2738.4734
$ws.Range("N132").Value = -18754.5005
$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N126").ClearContents()
$ws.Range("N137").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3960.45
$ws.Range("I7").Value = 3016.6667
$ws.Range("K7").Value = 3016.6667
$ws.Range("M7").Value = -2904.6667
$ws.Range("H22").Value = 41748.2
$ws.Range("I22").Value = 126425.25
$ws.Range("J22").Value = 1900.1765
$ws.Range("K22").Value = 126425.25
$ws.Range("L22").Value = 1900.1765
$ws.Range("M22").Value = -126130.25
$ws.Range("N22").Value = -2490.1765
$ws.Range("H27").Value = 41748.2
$ws.Range("I27").Value = 126425.25
$ws.Range("J27").Value = 1900.1765
$ws.Range("K27").Value = 126425.25
$ws.Range("L27").Value = 1900.1765
$ws.Range("M27").Value = -126318.25
$ws.Range("N27").Value = -2114.1765
$ws.Range("H29").Value = 14999.667
$ws.Range("J29").Value = 14999.667
$ws.Range("L29").Value = 14999.667
$ws.Range("N29").Value = -15589.667
$ws.Range("H40").Value = 6528.8613
$ws.Range("I40").Value = 5006.3687
$ws.Range("K40").Value = 5006.3687
$ws.Range("M40").Value = -4870.3687
$ws.Range("H46").Value = 2778
$ws.Range("I46").Value = 3460
$ws.Range("J46").Value = 2485.7144
$ws.Range("K46").Value = 3460
$ws.Range("L46").Value = 2485.7144
$ws.Range("M46").Value = -3272
$ws.Range("N46").Value = -2861.7144
$ws.Range("H55").Value = 381.6
$ws.Range("I55").Value = 334.9091
$ws.Range("J55").Value = 438.66666
$ws.Range("K55").Value = 334.9091
$ws.Range("L55").Value = 438.66666
$ws.Range("M55").Value = -161.9091
$ws.Range("N55").Value = -784.66666
$ws.Range("H93").Value = 8549305
$ws.Range("I93").Value = 15874700
$ws.Range("J93").Value = 3011.5
$ws.Range("K93").Value = 15874700
$ws.Range("L93").Value = 3011.5
$ws.Range("M93").Value = -15873452
$ws.Range("N93").Value = -5507.5
$ws.Range("H126").Value = 3960.45
$ws.Range("I126").Value = 3016.6667
$ws.Range("K126").Value = 9050.000100000001
$ws.Range("M126").Value = -6580.000100000001
$ws.Range("H132").Value = 14594.827
$ws.Range("I132").Value = 22242.584
$ws.Range("J132").Value = 9196.412
$ws.Range("K132").Value = 66727.75199999999
$ws.Range("L132").Value = 27589.236
$ws.Range("M132").Value = -64197.75199999999
$ws.Range("N132").Value = -32649.236
$ws.Range("H136").Value = 3359.484
$ws.Range("I136").Value = 1570.6364
$ws.Range("J136").Value = 7732.222
$ws.Range("K136").Value = 4711.9092
$ws.Range("L136").Value = 23196.666
$ws.Range("M136").Value = -2161.9092
$ws.Range("N136").Value = -28296.666

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 5966.6665
$ws.Range("I122").Value = 4542.857
$ws.Range("J122").Value = 7212.5
$ws.Range("K122").Value = 13628.571
$ws.Range("L122").Value = 21637.5
$ws.Range("M122").Value = -11178.571
$ws.Range("N122").Value = -26537.5
$ws.Range("H132").Value = 7248320.5
$ws.Range("I132").Value = 1472.5625
$ws.Range("J132").Value = 23812544
$ws.Range("K132").Value = 4417.6875
$ws.Range("L132").Value = 71437632
$ws.Range("M132").Value = -1887.6875
$ws.Range("N132").Value = -71442692
$ws.Range("H136").Value = 2038
$ws.Range("I136").Value = 643.1539
$ws.Range("J136").Value = 6571.25
$ws.Range("K136").Value = 1929.4617
$ws.Range("L136").Value = 19713.75
$ws.Range("M136").Value = 620.5382999999999
$ws.Range("N136").Value = -24813.75

